# Apply target sample size + selection fixes to the workbook
$wb = $excel.ActiveWorkbook

# --- Sheet "Vanilla" ---
$wsVanilla = $wb.Worksheets.Item("Vanilla")
$wsVanilla.Range("B2").Value = 120
$wsVanilla.Activate()
$wsVanilla.Range("B6").Select()

# --- Sheet "P8_Split_P6.2_only" ---
$wsP8 = $wb.Worksheets.Item("P8_Split_P6.2_only")
$wsP8.Range("B2").Value = 120
$wsP8.Range("B7").Value = 86
$wsP8.Activate()
$wsP8.Range("C8").Select()
